$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header cells O1:R1 (copy style+value from N1, then set values) ---
$ws.Range("N1").Copy($ws.Range("O1"))
$ws.Range("N1").Copy($ws.Range("P1"))
$ws.Range("N1").Copy($ws.Range("Q1"))
$ws.Range("N1").Copy($ws.Range("R1"))
$ws.Range("O1").Value = 13
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("R1").Value = 16

# --- Update existing value cells (rows 2-25) ---
# Row 2
$ws.Range("C2").Value = 0.9953323642099696
$ws.Range("D2").Value = 1.018278405596702
$ws.Range("E2").Value = 1.002489975646769
$ws.Range("I2").Value = 1.045792652206915
$ws.Range("J2").Value = 1.017653585419456
$ws.Range("K2").Value = 1.029488833613591
$ws.Range("L2").Value = 1.013914729811495
$ws.Range("N2").Value = 1.019098770237245
$ws.Range("Q2").Value = 1.02
$ws.Range("R2").Value = 1.031921327107745

# Row 3
$ws.Range("C3").Value = 0.999403245668435
$ws.Range("D3").Value = 1.021068537251389
$ws.Range("E3").Value = 1.0056943514573
$ws.Range("I3").Value = 1.04671233998829
$ws.Range("J3").Value = 1.019910192039134
$ws.Range("K3").Value = 1.031438756510595
$ws.Range("L3").Value = 1.016254886809278
$ws.Range("N3").Value = 1.021358581497164
$ws.Range("Q3").Value = 1.02
$ws.Range("R3").Value = 1.03329743003756

# Row 4
$ws.Range("C4").Value = 1.00198536701292
$ws.Range("D4").Value = 1.022838445406999
$ws.Range("E4").Value = 1.00773296021059
$ws.Range("I4").Value = 1.047280455023289
$ws.Range("J4").Value = 1.021338422885692
$ws.Range("K4").Value = 1.032668828278223
$ws.Range("L4").Value = 1.01773894262086
$ws.Range("N4").Value = 1.022788840595344
$ws.Range("Q4").Value = 1.02
$ws.Range("R4").Value = 1.034168115004029

# Row 5
$ws.Range("C5").Value = 1.003061441418702
$ws.Range("D5").Value = 1.023577820453777
$ws.Range("E5").Value = 1.008584263527116
$ws.Range("I5").Value = 1.047514751674555
$ws.Range("J5").Value = 1.021933945412402
$ws.Range("K5").Value = 1.033182090891077
$ws.Range("L5").Value = 1.018358146773567
$ws.Range("N5").Value = 1.02338520883235
$ws.Range("Q5").Value = 1.02
$ws.Range("R5").Value = 1.034538264701687

# Row 6
$ws.Range("C6").Value = 1.00324458037872
$ws.Range("D6").Value = 1.02370583735938
$ws.Range("E6").Value = 1.008729598203504
$ws.Range("I6").Value = 1.047555934591097
$ws.Range("J6").Value = 1.022036585165619
$ws.Range("K6").Value = 1.033272144668642
$ws.Range("L6").Value = 1.018464539083917
$ws.Range("N6").Value = 1.02348799434579
$ws.Range("Q6").Value = 1.02
$ws.Range("R6").Value = 1.0346107499119

# Row 7
$ws.Range("C7").Value = 1.002008175102083
$ws.Range("D7").Value = 1.022859910642975
$ws.Range("E7").Value = 1.00775198294596
$ws.Range("I7").Value = 1.047289485487823
$ws.Range("J7").Value = 1.021354578743397
$ws.Range("K7").Value = 1.032687136983715
$ws.Range("L7").Value = 1.017754769773701
$ws.Range("N7").Value = 1.02280501939622
$ws.Range("Q7").Value = 1.02
$ws.Range("R7").Value = 1.0342013245998

# Row 8
$ws.Range("C8").Value = 0.9967291369104473
$ws.Range("D8").Value = 1.019242590702799
$ws.Range("E8").Value = 1.003589389721203
$ws.Range("I8").Value = 1.046116101606695
$ws.Range("J8").Value = 1.018432647283239
$ws.Range("K8").Value = 1.030168067777888
$ws.Range("L8").Value = 1.014720958975126
$ws.Range("N8").Value = 1.019878938458239
$ws.Range("Q8").Value = 1.02
$ws.Range("R8").Value = 1.032424715824391

# Row 9
$ws.Range("C9").Value = 0.9869968329618741
$ws.Range("D9").Value = 1.012568776276837
$ws.Range("E9").Value = 0.9959641146622934
$ws.Range("I9").Value = 1.043828252365434
$ws.Range("J9").Value = 1.013017273022642
$ws.Range("K9").Value = 1.025461172862413
$ws.Range("L9").Value = 1.009123012437667
$ws.Range("N9").Value = 1.014455873745039
$ws.Range("Q9").Value = 1.02
$ws.Range("R9").Value = 1.029093355109321

# Row 10
$ws.Range("C10").Value = 0.9802021282206363
$ws.Range("D10").Value = 1.007917550451162
$ws.Range("E10").Value = 0.9906777109023882
$ws.Range("I10").Value = 1.042155476005237
$ws.Range("J10").Value = 1.009223868396925
$ws.Range("K10").Value = 1.022145307574688
$ws.Range("L10").Value = 1.005217347112521
$ws.Range("N10").Value = 1.010657082049644
$ws.Range("Q10").Value = 1.02
$ws.Range("R10").Value = 1.026765772854608

# Row 11
$ws.Range("C11").Value = 0.9771882528828518
$ws.Range("D11").Value = 1.005862716157315
$ws.Range("E11").Value = 0.9883435570670225
$ws.Range("I11").Value = 1.041399906925689
$ws.Range("J11").Value = 1.007541813468822
$ws.Range("K11").Value = 1.020674910728134
$ws.Range("L11").Value = 1.003488696399466
$ws.Range("N11").Value = 1.008972638410609
$ws.Range("Q11").Value = 1.02
$ws.Range("R11").Value = 1.025759354483176

# Row 12
$ws.Range("C12").Value = 0.976051352089391
$ws.Range("D12").Value = 1.005084648837552
$ws.Range("E12").Value = 0.9874638044664618
$ws.Range("I12").Value = 1.041110187733204
$ws.Range("J12").Value = 1.006905019119261
$ws.Range("K12").Value = 1.020115079060169
$ws.Range("L12").Value = 1.002835184778771
$ws.Range("N12").Value = 1.00833493974003
$ws.Range("Q12").Value = 1.02
$ws.Range("R12").Value = 1.025363525890429

# Row 13
$ws.Range("C13").Value = 0.9762948892534978
$ws.Range("D13").Value = 1.00525065370291
$ws.Range("E13").Value = 0.9876520469850897
$ws.Range("I13").Value = 1.041171953439938
$ws.Range("J13").Value = 1.007041082222261
$ws.Range("K13").Value = 1.020234247764593
$ws.Range("L13").Value = 1.002974856840854
$ws.Range("N13").Value = 1.008471196068243
$ws.Range("Q13").Value = 1.02
$ws.Range("R13").Value = 1.025445258957036

# Row 14
$ws.Range("C14").Value = 0.9770941453103992
$ws.Range("D14").Value = 1.005798024863787
$ws.Range("E14").Value = 0.9882706436197845
$ws.Range("I14").Value = 1.041375800764251
$ws.Range("J14").Value = 1.00748895423547
$ws.Range("K14").Value = 1.02062824755829
$ws.Range("L14").Value = 1.003434465282628
$ws.Range("N14").Value = 1.008919704111082
$ws.Range("Q14").Value = 1.02
$ws.Range("R14").Value = 1.025725279713998

# Row 15
$ws.Range("C15").Value = 0.9775867311727245
$ws.Range("D15").Value = 1.006136727104425
$ws.Range("E15").Value = 0.9886523686384591
$ws.Range("I15").Value = 1.04150192131846
$ws.Range("J15").Value = 1.007765659972493
$ws.Range("K15").Value = 1.020872544618367
$ws.Range("L15").Value = 1.003718367459326
$ws.Range("N15").Value = 1.00919680280199
$ws.Range("Q15").Value = 1.02
$ws.Range("R15").Value = 1.025903956497169

# Row 16
$ws.Range("C16").Value = 0.9804166092268158
$ws.Range("D16").Value = 1.008075248629968
$ws.Range("E16").Value = 0.990846484397619
$ws.Range("I16").Value = 1.04221618410119
$ws.Range("J16").Value = 1.00934999432659
$ws.Range("K16").Value = 1.022263962772251
$ws.Range("L16").Value = 1.005345828572965
$ws.Range("N16").Value = 1.010783387092596
$ws.Range("Q16").Value = 1.02
$ws.Range("R16").Value = 1.026890964539084

# Row 17
$ws.Range("C17").Value = 0.982169810127234
$ws.Range("D17").Value = 1.009276910474682
$ws.Range("E17").Value = 0.9922083273785325
$ws.Range("I17").Value = 1.042654154412991
$ws.Range("J17").Value = 1.010330850183811
$ws.Range("K17").Value = 1.023124213787429
$ws.Range("L17").Value = 1.006354366150974
$ws.Range("N17").Value = 1.01176563587764
$ws.Range("Q17").Value = 1.02
$ws.Range("R17").Value = 1.027501822362461

# Row 18
$ws.Range("C18").Value = 0.9831796679204077
$ws.Range("D18").Value = 1.009965739127566
$ws.Range("E18").Value = 0.9929929003587095
$ws.Range("I18").Value = 1.042902480578673
$ws.Range("J18").Value = 1.010893541591761
$ws.Range("K18").Value = 1.02361463996195
$ws.Range("L18").Value = 1.006933691669741
$ws.Range("N18").Value = 1.012329126371931
$ws.Range("Q18").Value = 1.02
$ws.Range("R18").Value = 1.027836802771735

# Row 19
$ws.Range("C19").Value = 0.9835266220239901
$ws.Range("D19").Value = 1.010205190190922
$ws.Range("E19").Value = 0.9932631788599344
$ws.Range("I19").Value = 1.042989296680132
$ws.Range("J19").Value = 1.011088389711266
$ws.Range("K19").Value = 1.023786452239589
$ws.Range("L19").Value = 1.007134046822424
$ws.Range("N19").Value = 1.012524251198115
$ws.Range("Q19").Value = 1.02
$ws.Range("R19").Value = 1.027964748810046

# Row 20
$ws.Range("C20").Value = 0.9819821198856855
$ws.Range("D20").Value = 1.009147967929571
$ws.Range("E20").Value = 0.9920623839222753
$ws.Range("I20").Value = 1.042607262005819
$ws.Range("J20").Value = 1.010225715253671
$ws.Range("K20").Value = 1.023031846257773
$ws.Range("L20").Value = 1.006246257629911
$ws.Range("N20").Value = 1.01166035164384
$ws.Range("Q20").Value = 1.02
$ws.Range("R20").Value = 1.027435162924037

# Row 21
$ws.Range("C21").Value = 0.9768642862746234
$ws.Range("D21").Value = 1.005644222004573
$ws.Range("E21").Value = 0.9880935014414409
$ws.Range("I21").Value = 1.041319548065873
$ws.Range("J21").Value = 1.007362192530772
$ws.Range("K21").Value = 1.020519446094555
$ws.Range("L21").Value = 1.003304012472623
$ws.Range("N21").Value = 1.008792762390224
$ws.Range("Q21").Value = 1.02
$ws.Range("R21").Value = 1.025662037324075

# Row 22
$ws.Range("C22").Value = 0.973566404827145
$ws.Range("D22").Value = 1.003384465527
$ws.Range("E22").Value = 0.9855435859382193
$ws.Range("I22").Value = 1.040472067222678
$ws.Range("J22").Value = 1.005512291848674
$ws.Range("K22").Value = 1.018889172336788
$ws.Range("L22").Value = 1.001406971959733
$ws.Range("N22").Value = 1.006940234636971
$ws.Range("Q22").Value = 1.02
$ws.Range("R22").Value = 1.024495664524969

# Row 23
$ws.Range("C23").Value = 0.9753157440744141
$ws.Range("D23").Value = 1.00457854672899
$ws.Range("E23").Value = 0.9868943758648139
$ws.Range("I23").Value = 1.040920244894559
$ws.Range("J23").Value = 1.006491326977697
$ws.Range("K23").Value = 1.019749121673701
$ws.Range("L23").Value = 1.002411057789112
$ws.Range("N23").Value = 1.007920660108174
$ws.Range("Q23").Value = 1.02
$ws.Range("R23").Value = 1.025093940233713

# Row 24
$ws.Range("C24").Value = 0.9820563568656405
$ws.Range("D24").Value = 1.009191542362616
$ws.Range("E24").Value = 0.9921185413778065
$ws.Range("I24").Value = 1.042621021727439
$ws.Range("J24").Value = 1.010263068703202
$ws.Range("K24").Value = 1.023059144284943
$ws.Range("L24").Value = 1.006285492532664
$ws.Range("N24").Value = 1.011697758139554
$ws.Range("Q24").Value = 1.02
$ws.Range("R24").Value = 1.027426745535829

# Row 25
$ws.Range("C25").Value = 0.9895769888701297
$ws.Range("D25").Value = 1.014345497311663
$ws.Range("E25").Value = 0.9979811177295721
$ws.Range("I25").Value = 1.044453772769078
$ws.Range("J25").Value = 1.014460293478956
$ws.Range("K25").Value = 1.026725125991639
$ws.Range("L25").Value = 1.010610855299135
$ws.Range("N25").Value = 1.015900943455919
$ws.Range("Q25").Value = 1.02
$ws.Range("R25").Value = 1.030015940391941
